# final edits to DoCalculations.py
# Update participant sample size (column B) from 10000 to 3000 for all
# simulation rows, and refresh the dependent results in columns E
# (participantAvgHullArea) and F (participantTime) that were recomputed
# after the sample-size change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  B = 3000; E = 7.002394562302324; F = 990.5599999999999 },
    @{ Row = 3;  B = 3000; E = 7.32185430412113;  F = 1412.5895 },
    @{ Row = 4;  B = 3000; E = 7.85893764337613;  F = 642.5690000000001 },
    @{ Row = 5;  B = 3000; E = 5.283706020567705; F = 548.7927000000002 },
    @{ Row = 6;  B = 3000; E = 5.597256745794178; F = 475.0768 },
    @{ Row = 7;  B = 3000; E = 5.157439762774044; F = 650.6723 },
    @{ Row = 8;  B = 3000; E = 6.006502478787697; F = 993.6983000000001 },
    @{ Row = 9;  B = 3000; E = 5.204135757628497; F = 507.3191999999999 },
    @{ Row = 10; B = 3000; E = 8.600099510621336; F = 522.0576000000001 },
    @{ Row = 11; B = 3000; E = 7.590580015361675; F = 1051.0824 },
    @{ Row = 12; B = 3000; E = 6.640301149951249; F = 1486.608300000001 },
    @{ Row = 13; B = 3000; E = 9.660098499283437; F = 117.4134 },
    @{ Row = 14; B = 3000; E = 5.606104369014449; F = 1271.6492 },
    @{ Row = 15; B = 3000; E = 4.065004797632901; F = 810.3334000000001 },
    @{ Row = 16; B = 3000; E = 4.740030767301473; F = 554.1023 },
    @{ Row = 17; B = 3000; E = 5.148374528378501; F = 572.0085999999999 },
    @{ Row = 18; B = 3000; E = 5.300448792010487; F = 805.2714999999998 },
    @{ Row = 19; B = 3000; E = 3.302175817650669; F = 739.8534999999999 },
    @{ Row = 20; B = 3000; E = 5.32056935755115;  F = 993.7463 },
    @{ Row = 21; B = 3000; E = 4.855349531467875; F = 541.4542000000001 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B   # column B - participant
    $ws.Cells.Item($u.Row, 5).Value = $u.E   # column E - participantAvgHullArea
    $ws.Cells.Item($u.Row, 6).Value = $u.F   # column F - participantTime
}
